$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.105.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.882.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.30%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.882.76"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.35%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.73%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.38"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.19"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.68%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.529.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.884.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.176.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.50%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.01%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.90%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.27%  "

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.028.40"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.20%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.20"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.825.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.12%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.42%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.93"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.55%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "448.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.74%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.83%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.31%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.59"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.22%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.50"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +14.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.827.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.98%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.18%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0359"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.37%  "
